$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A5 loses its old text "awd" but keeps its style (s=1), left an empty numeric cell.
$ws.Range("A5").Value = $null

# C5 gets the same (centered) style as A1 / C2, with new text "Direct Flight".
$ws.Range("A1").Copy()
$ws.Range("C5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C5").Value = "Direct Flight"
$excel.CutCopyMode = 0

# D5 must hold "5591" as literal TEXT, not be auto-coerced to a number, and
# must NOT pick up a distinct number-format style (stays default style 0).
# Build it as a text formula, then flatten to a plain value in place.
$ws.Range("D5").Formula = "=""5591"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0

$ws.Range("E5").Value = "09AUG22-11AUG22"
$ws.Range("F5").Value = ".2.4..."

# G5 / H5 are time-of-day values formatted h:mm:ss (numFmtId 21) - new shared style.
$ws.Range("G5").Value = 0.4791666666666667
$ws.Range("G5").NumberFormat = "h:mm:ss"

$ws.Range("H5").Value = 0.5590277777777778
$ws.Range("H5").NumberFormat = "h:mm:ss"

$ws.Range("I5").Value = 0
